# Auto-generated edit script updating H:N (market price / profit) columns
# across the Leve-profit worksheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 6585.4287
$ws.Range("I4").Value = 3024.5
$ws.Range("J4").Value = 11333.333
$ws.Range("K4").Value = 3024.5
$ws.Range("L4").Value = 11333.333
$ws.Range("M4").Value = -2910.5
$ws.Range("N4").Value = -11561.333

$ws.Range("H86").Value = 650762.3
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 650762.3
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H96").Value = 8762.538
$ws.Range("I96").Value = 11780
$ws.Range("J96").Value = 1973.25
$ws.Range("K96").Value = 35340
$ws.Range("L96").Value = 5919.75
$ws.Range("M96").Value = -33967
$ws.Range("N96").Value = -8665.75

$ws.Range("H115").Value = 799
$ws.Range("I115").Value = 799
$ws.Range("K115").Value = 2397
$ws.Range("M115").Value = -830

$ws.Range("H127").Value = 2037
$ws.Range("I127").Value = 1959.8334
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 5879.5002
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = -919.5002000000004
$ws.Range("N127").Value = -17420

$ws.Range("H137").Value = 1000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws.Range("H138").Value = 1501.1578
$ws.Range("I138").Value = 1529.6666
$ws.Range("J138").Value = 988
$ws.Range("K138").Value = 4588.9998
$ws.Range("L138").Value = 2964
$ws.Range("M138").Value = 551.0002000000004
$ws.Range("N138").Value = -13244


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 115.14286

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H32").Value = 4424729.5
$ws.Range("I32").Value = 4121756.2
$ws.Range("K32").Value = 4121756.2
$ws.Range("M32").Value = -4121469.2

$ws.Range("H86").Value = 27000
$ws.Range("J86").Value = 27000
$ws.Range("L86").Value = 27000
$ws.Range("N86").Value = -29372

$ws.Range("H89").Value = 27000
$ws.Range("J89").Value = 27000
$ws.Range("L89").Value = 81000
$ws.Range("N89").Value = -92856

$ws.Range("H109").Value = 89999.5
$ws.Range("J109").Value = 89999.5
$ws.Range("L109").Value = 89999.5
$ws.Range("N109").Value = -92773.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 115.14286

$ws.Range("H86").Value = 1244.4375
$ws.Range("I86").Value = 1016.1111
$ws.Range("K86").Value = 1016.1111
$ws.Range("M86").Value = 106.8889

$ws.Range("H89").Value = 1244.4375
$ws.Range("I89").Value = 1016.1111
$ws.Range("K89").Value = 5080.555499999999
$ws.Range("M89").Value = 535.4445000000005


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H7").Value = 1767.909
$ws.Range("I7").Value = 492.7143
$ws.Range("J7").Value = 3999.5
$ws.Range("K7").Value = 492.7143
$ws.Range("L7").Value = 3999.5
$ws.Range("M7").Value = -379.7143
$ws.Range("N7").Value = -4225.5

$ws.Range("H58").Value = 997
$ws.Range("I58").Value = 997
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 997
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -794
$ws.Range("N58").ClearContents()

$ws.Range("H86").Value = 7949.75
$ws.Range("I86").Value = 6599
$ws.Range("K86").Value = 6599
$ws.Range("M86").Value = -5476

$ws.Range("H89").Value = 7949.75
$ws.Range("I89").Value = 6599
$ws.Range("K89").Value = 32995
$ws.Range("M89").Value = -27379

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 5902.643
$ws.Range("J132").Value = 4372
$ws.Range("L132").Value = 13116
$ws.Range("N132").Value = -18176

$ws.Range("H136").Value = 997
$ws.Range("I136").Value = 997
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2991
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -441
$ws.Range("N136").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 13990.363
$ws.Range("I121").Value = 18733.334
$ws.Range("K121").Value = 56200.00199999999
$ws.Range("M121").Value = -54890.00199999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 45705
$ws.Range("I5").Value = 78980
$ws.Range("K5").Value = 78980
$ws.Range("M5").Value = -78868

$ws.Range("H29").Value = 35000
$ws.Range("J29").Value = 35000
$ws.Range("L29").Value = 35000
$ws.Range("N29").Value = -35580

$ws.Range("H80").Value = 3749.5
$ws.Range("I80").Value = 3749.5
$ws.Range("K80").Value = 3749.5
$ws.Range("M80").Value = -2751.5

$ws.Range("H83").Value = 3749.5
$ws.Range("I83").Value = 3749.5
$ws.Range("K83").Value = 18747.5
$ws.Range("M83").Value = -13755.5

$ws.Range("H86").Value = 43529
$ws.Range("J86").Value = 43529
$ws.Range("L86").Value = 43529
$ws.Range("N86").Value = -45901

$ws.Range("H89").Value = 43529
$ws.Range("J89").Value = 43529
$ws.Range("L89").Value = 130587
$ws.Range("N89").Value = -142443


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3191.3333
$ws.Range("I22").Value = 3480
$ws.Range("J22").Value = 1748
$ws.Range("K22").Value = 3480
$ws.Range("L22").Value = 1748
$ws.Range("M22").Value = -3185
$ws.Range("N22").Value = -2338

$ws.Range("H27").Value = 3191.3333
$ws.Range("I27").Value = 3480
$ws.Range("J27").Value = 1748
$ws.Range("K27").Value = 3480
$ws.Range("L27").Value = 1748
$ws.Range("M27").Value = -3373
$ws.Range("N27").Value = -1962

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H132").Value = 3305.8572
$ws.Range("I132").Value = 3108.8
$ws.Range("K132").Value = 9326.400000000001
$ws.Range("M132").Value = -6796.400000000001

